$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update headers (row 1)
$ws.Range("A1").Value = "Vehicle ID"
$ws.Range("B1").Value = "Vehicle Name"
$ws.Range("C1").Value = "Fuel Type"
$ws.Range("D1").Value = "Mileage Unit"
$ws.Range("E1").Value = "Mileage Value"
$ws.Range("F1").Value = "Remarks"

# Update data row 2
$ws.Range("A2").Value = "V101"
$ws.Range("B2").Value = "Bus"
$ws.Range("C2").Value = "disel:f102"
$ws.Range("D2").Value = "mpg"
$ws.Range("E2").Value = 15.0
$ws.Range("F2").Value = "This is a commuter bus"

# Clear old row 3 (no longer present in the new layout)
$ws.Range("A3:F3").Clear()
